$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista")

$ws.Range("A9").Value = "Empresa 4"
$ws.Range("B9").Value = "Azúcar x1kg"
$ws.Range("C9").Value = 800

$ws.Range("A10").Value = "Empresa 4"
$ws.Range("B10").Value = "Sal x1kg"
$ws.Range("C10").Value = 900

$ws.Range("C10").Select()
